$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4293.6875
$ws.Range("I116").Value = 1956
$ws.Range("K116").Value = 1956
$ws.Range("M116").Value = 1486
$ws.Range("H127").Value = 1312.7858
$ws.Range("I127").Value = 519.9
$ws.Range("J127").Value = 3295
$ws.Range("K127").Value = 1559.7
$ws.Range("L127").Value = 9885
$ws.Range("M127").Value = 3400.3
$ws.Range("N127").Value = -19805
$ws.Range("H138").Value = 2077.652
$ws.Range("I138").Value = 1689.4546
$ws.Range("J138").Value = 2199.6572
$ws.Range("K138").Value = 5068.3638
$ws.Range("L138").Value = 6598.971600000001
$ws.Range("M138").Value = 71.63619999999992
$ws.Range("N138").Value = -16878.9716
$ws.Range("H141").Value = 2638.3333
$ws.Range("I141").Value = 2405.625
$ws.Range("J141").Value = 4500
$ws.Range("K141").Value = 7216.875
$ws.Range("L141").Value = 13500
$ws.Range("M141").Value = -2036.875
$ws.Range("N141").Value = -23860

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1468.2972
$ws.Range("I2").Value = 1461.1538
$ws.Range("J2").Value = 1485.1818
$ws.Range("K2").Value = 1461.1538
$ws.Range("L2").Value = 1485.1818
$ws.Range("M2").Value = -1348.1538
$ws.Range("N2").Value = -1711.1818
$ws.Range("H32").Value = 2965.5862
$ws.Range("I32").Value = 2290.4727
$ws.Range("J32").Value = 15342.667
$ws.Range("K32").Value = 2290.4727
$ws.Range("L32").Value = 15342.667
$ws.Range("M32").Value = -2003.4727
$ws.Range("N32").Value = -15916.667
$ws.Range("H60").Value = 16633.334
$ws.Range("I60").Value = 11900
$ws.Range("J60").Value = 19000
$ws.Range("K60").Value = 11900
$ws.Range("L60").Value = 19000
$ws.Range("M60").Value = -11167
$ws.Range("N60").Value = -20466
$ws.Range("I63").Value = 31249998
$ws.Range("K63").Value = 31249998
$ws.Range("M63").Value = -31249312
$ws.Range("I66").Value = 31249998
$ws.Range("K66").Value = 156249990
$ws.Range("M66").Value = -156246558
$ws.Range("H116").Value = 1468.2972
$ws.Range("I116").Value = 1461.1538
$ws.Range("J116").Value = 1485.1818
$ws.Range("K116").Value = 1461.1538
$ws.Range("L116").Value = 1485.1818
$ws.Range("M116").Value = 832.8462
$ws.Range("N116").Value = -6073.1818

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1468.2972
$ws.Range("I3").Value = 1461.1538
$ws.Range("J3").Value = 1485.1818
$ws.Range("K3").Value = 1461.1538
$ws.Range("L3").Value = 1485.1818
$ws.Range("M3").Value = -1347.1538
$ws.Range("N3").Value = -1713.1818
$ws.Range("H20").Value = 2623.5
$ws.Range("I20").Value = 3548
$ws.Range("J20").Value = 959.4
$ws.Range("K20").Value = 3548
$ws.Range("L20").Value = 959.4
$ws.Range("M20").Value = -3301
$ws.Range("N20").Value = -1453.4
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H134").Value = 4448.923
$ws.Range("I134").Value = 4697.6665
$ws.Range("J134").Value = 1464
$ws.Range("K134").Value = 14092.9995
$ws.Range("L134").Value = 4392
$ws.Range("M134").Value = -11557.9995
$ws.Range("N134").Value = -9462
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10010.909
$ws.Range("I31").Value = 14762.167
$ws.Range("J31").Value = 4309.4
$ws.Range("K31").Value = 14762.167
$ws.Range("L31").Value = 4309.4
$ws.Range("M31").Value = -14467.167
$ws.Range("N31").Value = -4899.4
$ws.Range("H34").Value = 10010.909
$ws.Range("I34").Value = 14762.167
$ws.Range("J34").Value = 4309.4
$ws.Range("K34").Value = 14762.167
$ws.Range("L34").Value = 4309.4
$ws.Range("M34").Value = -14560.167
$ws.Range("N34").Value = -4713.4
$ws.Range("H52").Value = 33800
$ws.Range("J52").Value = 33800
$ws.Range("L52").Value = 33800
$ws.Range("N52").Value = -34388
$ws.Range("H58").Value = 24147.545
$ws.Range("I58").Value = 1740.5
$ws.Range("J58").Value = 51036
$ws.Range("K58").Value = 1740.5
$ws.Range("L58").Value = 51036
$ws.Range("M58").Value = -1537.5
$ws.Range("N58").Value = -51442
$ws.Range("H105").Value = 1988.4286
$ws.Range("I105").Value = 1983.8
$ws.Range("K105").Value = 1983.8
$ws.Range("M105").Value = -236.8
$ws.Range("H134").Value = 1309.6945
$ws.Range("I134").Value = 970.04
$ws.Range("J134").Value = 2081.6365
$ws.Range("K134").Value = 2910.12
$ws.Range("L134").Value = 6244.9095
$ws.Range("M134").Value = -375.1199999999999
$ws.Range("N134").Value = -11314.9095
$ws.Range("H136").Value = 24147.545
$ws.Range("I136").Value = 1740.5
$ws.Range("J136").Value = 51036
$ws.Range("K136").Value = 5221.5
$ws.Range("L136").Value = 153108
$ws.Range("M136").Value = -2671.5
$ws.Range("N136").Value = -158208

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 775.23
$ws.Range("J131").Value = 775.23
$ws.Range("L131").Value = 2325.69
$ws.Range("N131").Value = -12405.69

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3182.1785
$ws.Range("I80").Value = 2876.25
$ws.Range("J80").Value = 3411.625
$ws.Range("K80").Value = 2876.25
$ws.Range("L80").Value = 3411.625
$ws.Range("M80").Value = -1878.25
$ws.Range("N80").Value = -5407.625
$ws.Range("H83").Value = 3182.1785
$ws.Range("I83").Value = 2876.25
$ws.Range("J83").Value = 3411.625
$ws.Range("K83").Value = 14381.25
$ws.Range("L83").Value = 17058.125
$ws.Range("M83").Value = -9389.25
$ws.Range("N83").Value = -27042.125
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 2738.7222
$ws.Range("I113").Value = 2327
$ws.Range("K113").Value = 2327
$ws.Range("M113").Value = -157
$ws.Range("H132").Value = 27505.455
$ws.Range("I132").Value = 5439
$ws.Range("J132").Value = 66121.75
$ws.Range("K132").Value = 16317
$ws.Range("L132").Value = 198365.25
$ws.Range("M132").Value = -13787
$ws.Range("N132").Value = -203425.25

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3694.32
$ws.Range("I61").Value = 2362
$ws.Range("J61").Value = 6525.5
$ws.Range("K61").Value = 2362
$ws.Range("L61").Value = 6525.5
$ws.Range("M61").Value = -2160
$ws.Range("N61").Value = -6929.5
$ws.Range("H113").Value = 3694.32
$ws.Range("I113").Value = 2362
$ws.Range("J113").Value = 6525.5
$ws.Range("K113").Value = 2362
$ws.Range("L113").Value = 6525.5
$ws.Range("M113").Value = -192
$ws.Range("N113").Value = -10865.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5105.75
$ws.Range("I62").Value = 4885
$ws.Range("J62").Value = 5149.9
$ws.Range("K62").Value = 4885
$ws.Range("L62").Value = 5149.9
$ws.Range("M62").Value = -4261
$ws.Range("N62").Value = -6397.9
$ws.Range("H65").Value = 5105.75
$ws.Range("I65").Value = 4885
$ws.Range("J65").Value = 5149.9
$ws.Range("K65").Value = 24425
$ws.Range("L65").Value = 25749.5
$ws.Range("M65").Value = -21305
$ws.Range("N65").Value = -31989.5
